$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the duplicate "Play Aliens and Pyramid..." bold paragraph
#    that sits right before the final italic meta-description paragraph,
#    and update that meta-description text to the new "Prompt: ..." copy.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Play Aliens and Pyramid Free Slot - A Unique and Thrilling Gaming Experience`r") {
        if ($i -gt 1) {
            $p.Range.Delete()
        }
    }
}

$oldDesc = "Read our review of Aliens and Pyramid slot game. Play for free and enjoy a thrilling and exciting gaming experience with unique graphics and theme."
$newPrompt = "Prompt: Create a feature image for Aliens and Pyramids, the online slot game. The image should be in cartoon style and feature a happy Maya warrior with glasses. The image should be eye-catching and entice players to try out the game. The Maya warrior in the image should hold a golden key to depict the theme of the game, which is based on aliens and ancient Egypt. Use bright colors and be creative in designing the image."

$d.Content.Find.Execute($oldDesc, $true, $false, $false, $false, $false, $true, 1, $false, $newPrompt, 2)

# ------------------------------------------------------------------
# 2) Insert a new "Meta description" paragraph right after the H1 title.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(2)
$newPara.Style = "Normal"

$metaLabel = "Meta description"
$metaRest = ": Read our review of Aliens and Pyramid slot game. Play for free and enjoy a thrilling and exciting gaming experience with unique graphics and theme."

$insertionPoint = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$insertionPoint.InsertAfter($metaLabel + $metaRest)

$boldRange = $d.Range($newPara.Range.Start, $newPara.Range.Start + $metaLabel.Length)
$boldRange.Bold = 1
